$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Grade values (shared string "B"/"A" -> "C") for the rows touched
# by the upload: No. 11, 29, 34, 43, 44, 54, 64, 65, 68, 76, 87
$ws.Range("B12").Value = "C"
$ws.Range("B30").Value = "C"
$ws.Range("B35").Value = "C"
$ws.Range("B44").Value = "C"
$ws.Range("B45").Value = "C"
$ws.Range("B55").Value = "C"
$ws.Range("B65").Value = "C"
$ws.Range("B66").Value = "C"
$ws.Range("B69").Value = "C"
$ws.Range("B77").Value = "C"
$ws.Range("B88").Value = "C"

# Move the view / selection to match the author's last on-screen position
$ws.Range("E74").Select()
